$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 data: replace the Parent/BIO/Video (parents.com article) row
# with the instyle/TAXONOMY row, and drop the viewType (column D) value.
$ws.Range("A3").Value = "https://www-instyle-ddm-staging.a-ue1.instyle.com/"
$ws.Range("B3").Value = "instyle"
$ws.Range("C3").Value = "TAXONOMY"
$ws.Range("D3").Value = ""

# Remove the hyperlink that was attached to A3 (the old parents.com article URL).
$ws.Hyperlinks.Delete()

# Update the active selection to D3, matching the saved view state.
$ws.Range("D3").Select()
